# Swap the fynd ("observation") data between row 2 and row 4.
# Columns that carry distinct data for these two records:
#   A=1 Id, B=2 Taxonsorteringsordning, E=5 TaxonId, F=6 Artnamn,
#   G=7 Vetenskapligt namn, H=8 Auktor, Q=17 Ost, R=18 Nord
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @(1, 2, 5, 6, 7, 8, 17, 18)

foreach ($c in $cols) {
    $v2 = $ws.Cells.Item(2, $c).Value2
    $v4 = $ws.Cells.Item(4, $c).Value2
    $ws.Cells.Item(2, $c).Value2 = $v4
    $ws.Cells.Item(4, $c).Value2 = $v2
}
